# Update ID labels on the "Feuil1" sheet:
#   get_started  -> phrase_presentation  (rows 2 and 3, column A)
#   out_of_scope -> phrase_hors_sujet    (row 4, column A)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("A2").Value = "phrase_presentation"
$ws.Range("A3").Value = "phrase_presentation"
$ws.Range("A4").Value = "phrase_hors_sujet"

# Move the active selection to A3, matching the saved state of the file.
$ws.Activate()
$ws.Range("A3").Select()
